$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 47: no content changes (date 45175, B47/C47/D47 unchanged). Nothing to
# edit here - its style index shift is a pure side effect of the new style
# entry and requires no explicit action.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 48 gains a new, empty "A48" cell so that A47:A49 can be merged as one
# date block. Give it the same look as the other date cells (centered,
# numFmtId 14) by copying the style from A47, then blanking the value.
# ---------------------------------------------------------------------------
$ws.Range("A47").Copy($ws.Range("A48"))
$ws.Range("A48").Value = ""

# ---------------------------------------------------------------------------
# Row 49 (new): same date block as 47/48 (empty A49), a new time range in B
# and a new task description in C.
# ---------------------------------------------------------------------------
$ws.Range("A47").Copy($ws.Range("A49"))
$ws.Range("A49").Value = ""
$ws.Range("B49").Value = "13:30 - 20:00"
$ws.Range("C47").Copy($ws.Range("C49"))
$ws.Range("C49").Value = 'Lecture de la documentation et programmation du module RFID "CHILLI"'

# Merge the three-row date block for column A.
$ws.Range("A47:A49").Merge()

# ---------------------------------------------------------------------------
# Row 50 (new): next day (07.09.2023 / serial 45176), a recurring time
# range, and a new task description (styled like the wrapped/merged C45:C46
# block).
# ---------------------------------------------------------------------------
$ws.Range("A47").Copy($ws.Range("A50"))
$ws.Range("A50").Value = 45176
$ws.Range("B50").Value = "08:30 - 11:30"
$ws.Range("C45").Copy($ws.Range("C50"))
$ws.Range("C50").Value = 'Programmation du module RFID "CHILLI" avec utilisation des librairies du fabricant'

# ---------------------------------------------------------------------------
# Row 51 (new): closes out the day - empty date cell (plain centered, no
# date format this time), a recurring time label, and an empty, merged
# C-cell to close the C50:C51 block.
# ---------------------------------------------------------------------------
$ws.Range("B22").Copy($ws.Range("B51"))
$ws.Range("B51").Value = "13:00 - 16:00"
$ws.Range("C45").Copy($ws.Range("C51"))
$ws.Range("C51").Value = ""

# Merge the new date block and the wrapped task-description block first...
$ws.Range("A50:A51").Merge()
$ws.Range("C50:C51").Merge()

# ...then give A51 its own (non-date) look: merging pulls A50's date format
# onto the whole range, so restyle the bottom cell afterwards.
$ws.Range("B1").Copy($ws.Range("A51"))
$ws.Range("A51").Value = ""

$ws.Range("D50").Select()
